# Created Success Page, Excel updation, vectors
# Appends 5 new payment-log rows (26-30) to the "Payments" sheet, mirroring
# the rows written by the app's Razorpay webhook logger.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns that hold values which *look* numeric/date/time but must be kept
# as literal text (matching how the original rows were authored) need their
# number format forced to Text ("@") before the value is assigned, otherwise
# Excel's automatic type inference would turn them into a date serial /
# number and silently mangle things like a leading "+" or a leading "0".
$textCols = @("A", "B", "C", "E", "J")
foreach ($col in $textCols) {
    $rangeAddr = $col + "26:" + $col + "30"
    $ws.Range($rangeAddr).NumberFormat = "@"
}

$rows = @(
    @("1/1/2025", "1:11:58 pm", "010125131158", "order_Pe6HZcMA3jtKOF", "21", "21B81A05V9", "SAMRATH REDDY", "CSE", "E", "+917999999990", "fomowog893@nozamas.com", "CollegeFee", "I",  "", 120000, "wallet", "Verification in progress..."),
    @("1/1/2025", "1:15:15 pm", "010125131515", "order_Pe6KAtrGxllVGc", "21", "21B81A05V9", "SAMRATH REDDY", "CSE", "E", "+917999999990", "fomowog893@nozamas.com", "CollegeFee", "I",  "", 120000, "wallet", "Verification in progress..."),
    @("1/1/2025", "3:10:43 pm", "010125151043", "order_Pe8IxBmB7VvSbS", "21", "21B81A05V9", "SAMRATH REDDY", "CSE", "E", "+917981455290", "samrathreddy04@gmail.com", "CollegeFee", "IV", "", 120000, "wallet", "Verification in progress..."),
    @("1/1/2025", "3:30:19 pm", "010125153019", "order_Pe8dccfdSEBR2X", "21", "21B81A05V9", "SAMRATH REDDY", "CSE", "E", "+917981455290", "samrathreddy04@gmail.com", "CollegeFee", "IV", "", 120000, "wallet", "Rejected"),
    @("1/1/2025", "3:35:52 pm", "010125153552", "order_Pe8jN3nSLdJf6t", "21", "21B81A05V9", "SAMRATH REDDY", "CSE", "E", "+917981455290", "samrathreddy04@gmail.com", "CollegeFee", "I",  "", 120000, "wallet", "Verified")
)

$startRow = 26
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $values = $rows[$i]
    for ($c = 0; $c -lt $values.Count; $c++) {
        $val = $values[$c]
        # Column N (index 13) is blank in every appended row (FeeSem is
        # unused for these transactions) - skip it so the cell is left
        # empty instead of writing a zero-length string.
        if ($c -eq 13) {
            continue
        }
        $ws.Cells.Item($r, $c + 1).Value = $val
    }
}

Write-Host "Appended rows 26-30 to sheet '$($ws.Name)'"
